$wb = $excel.ActiveWorkbook

# Rename existing Sheet1 to ValidLogin
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

# Add a new worksheet for InvalidLogin, placed after ValidLogin
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

# Fill ValidLogin sheet data
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("C1").Value = "FailMsg"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
$ws1.Range("C2").Value = "Home Page is not displayed"

# Fill InvalidLogin sheet data
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "FailMsg"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "test"
$ws2.Range("C2").Value = "Home Page is not displayed"

# Set column width on ValidLogin (bestFit width 6.57 ~ autofit)
$ws1.Columns.Item(1).AutoFit() | Out-Null

# Selections
$ws1.Range("A1:E2").Select()
$ws2.Range("F10").Select()

# Activate InvalidLogin sheet (tabSelected / activeTab = 1)
$ws2.Activate()
